$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 3327.2  # H19: 4010 -> 3327.2
$ws.Cells.Item(19, 9).Value = 3976.4  # I19: 6230 -> 3976.4
$ws.Cells.Item(19, 11).Value = 3976.4  # K19: 6230 -> 3976.4
$ws.Cells.Item(19, 13).Value = -3801.4  # M19: -6055 -> -3801.4

$ws.Cells.Item(28, 8).Value = 850.2353000000001  # H28: 909.1875 -> 850.2353000000001
$ws.Cells.Item(28, 9).Value = 876.13336  # I28: 994.1539 -> 876.13336
$ws.Cells.Item(28, 10).Value = 656  # J28: 541 -> 656
$ws.Cells.Item(28, 11).Value = 876.13336  # K28: 994.1539 -> 876.13336
$ws.Cells.Item(28, 12).Value = 656  # L28: 541 -> 656
$ws.Cells.Item(28, 13).Value = -391.13336  # M28: -509.1539 -> -391.13336
$ws.Cells.Item(28, 14).Value = -1626  # N28: -1511 -> -1626

$ws.Cells.Item(40, 8).Value = 4450  # H40: 5666.6665 -> 4450
$ws.Cells.Item(40, 9).Value = 3400  # I40: 0 -> 3400
$ws.Cells.Item(40, 10).Value = 5500  # J40: 5666.6665 -> 5500
$ws.Cells.Item(40, 11).Value = 3400  # K40: 0 -> 3400
$ws.Cells.Item(40, 12).Value = 5500  # L40: 5666.6665 -> 5500
$ws.Cells.Item(40, 13).Value = -3225  # M40: add -3225
$ws.Cells.Item(40, 14).Value = -5850  # N40: -6016.6665 -> -5850

$ws.Cells.Item(62, 8).Value = 71791.8  # H62: 88540.164 -> 71791.8
$ws.Cells.Item(62, 9).Value = 95307.73  # I62: 104448.2 -> 95307.73
$ws.Cells.Item(62, 10).Value = 7123  # J62: 9000 -> 7123
$ws.Cells.Item(62, 11).Value = 95307.73  # K62: 104448.2 -> 95307.73
$ws.Cells.Item(62, 12).Value = 7123  # L62: 9000 -> 7123
$ws.Cells.Item(62, 13).Value = -94683.73  # M62: -103824.2 -> -94683.73
$ws.Cells.Item(62, 14).Value = -8371  # N62: -10248 -> -8371

$ws.Cells.Item(65, 8).Value = 71791.8  # H65: 88540.164 -> 71791.8
$ws.Cells.Item(65, 9).Value = 95307.73  # I65: 104448.2 -> 95307.73
$ws.Cells.Item(65, 10).Value = 7123  # J65: 9000 -> 7123
$ws.Cells.Item(65, 11).Value = 476538.65  # K65: 522241 -> 476538.65
$ws.Cells.Item(65, 12).Value = 35615  # L65: 45000 -> 35615
$ws.Cells.Item(65, 13).Value = -473418.65  # M65: -519121 -> -473418.65
$ws.Cells.Item(65, 14).Value = -41855  # N65: -51240 -> -41855

$ws.Cells.Item(112, 8).Value = 2098.4  # H112: 2113.2 -> 2098.4
$ws.Cells.Item(112, 10).Value = 2232.125  # J112: 2250.625 -> 2232.125
$ws.Cells.Item(112, 12).Value = 6696.375  # L112: 6751.875 -> 6696.375
$ws.Cells.Item(112, 14).Value = -8912.375  # N112: -8967.875 -> -8912.375

$ws.Cells.Item(132, 8).Value = 4161.1113  # H132: 4265.6855 -> 4161.1113
$ws.Cells.Item(132, 9).Value = 1879.9565  # I132: 1942.6364 -> 1879.9565
$ws.Cells.Item(132, 11).Value = 5639.8695  # K132: 5827.9092 -> 5639.8695
$ws.Cells.Item(132, 13).Value = -3109.8695  # M132: -3297.9092 -> -3109.8695

$ws.Cells.Item(135, 8).Value = 855.5357  # H135: 34150.867 -> 855.5357
$ws.Cells.Item(135, 9).Value = 782.2  # I135: 774.0769 -> 782.2
$ws.Cells.Item(135, 10).Value = 1466.6666  # J135: 251100 -> 1466.6666
$ws.Cells.Item(135, 11).Value = 7039.8  # K135: 6966.6921 -> 7039.8
$ws.Cells.Item(135, 12).Value = 13199.9994  # L135: 2259900 -> 13199.9994
$ws.Cells.Item(135, 13).Value = -4504.8  # M135: -4431.6921 -> -4504.8
$ws.Cells.Item(135, 14).Value = -18269.9994  # N135: -2264970 -> -18269.9994

$ws.Cells.Item(137, 8).Value = 24490.582  # H137: 26917.82 -> 24490.582
$ws.Cells.Item(137, 9).Value = 1235.125  # I137: 1265.2903 -> 1235.125
$ws.Cells.Item(137, 10).Value = 92142.82000000001  # J137: 126321.375 -> 92142.82000000001
$ws.Cells.Item(137, 11).Value = 3705.375  # K137: 3795.8709 -> 3705.375
$ws.Cells.Item(137, 12).Value = 276428.46  # L137: 378964.125 -> 276428.46
$ws.Cells.Item(137, 13).Value = -1155.375  # M137: -1245.8709 -> -1155.375
$ws.Cells.Item(137, 14).Value = -281528.46  # N137: -384064.125 -> -281528.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 6000  # H25: 0 -> 6000
$ws.Cells.Item(25, 10).Value = 6000  # J25: 0 -> 6000
$ws.Cells.Item(25, 12).Value = 6000  # L25: 0 -> 6000
$ws.Cells.Item(25, 13).Value = -6804  # M25: add -6804

$ws.Cells.Item(31, 8).Value = 8161.8887  # H31: 10508.286 -> 8161.8887
$ws.Cells.Item(31, 9).Value = 5367.25  # I31: 7173.1665 -> 5367.25
$ws.Cells.Item(31, 11).Value = 5367.25  # K31: 7173.1665 -> 5367.25
$ws.Cells.Item(31, 13).Value = -5073.25  # M31: -6879.1665 -> -5073.25

$ws.Cells.Item(32, 8).Value = 39998.723  # H32: 41438.44 -> 39998.723
$ws.Cells.Item(32, 9).Value = 22102.896  # I32: 22570.361 -> 22102.896
$ws.Cells.Item(32, 10).Value = 183165.33  # J32: 218798.4 -> 183165.33
$ws.Cells.Item(32, 11).Value = 22102.896  # K32: 22570.361 -> 22102.896
$ws.Cells.Item(32, 12).Value = 183165.33  # L32: 218798.4 -> 183165.33
$ws.Cells.Item(32, 13).Value = -21815.896  # M32: -22283.361 -> -21815.896
$ws.Cells.Item(32, 14).Value = -183739.33  # N32: -219372.4 -> -183739.33

$ws.Cells.Item(35, 8).Value = 2456.25  # H35: 0 -> 2456.25
$ws.Cells.Item(35, 9).Value = 2456.25  # I35: 0 -> 2456.25
$ws.Cells.Item(35, 11).Value = 2456.25  # K35: 0 -> 2456.25
$ws.Cells.Item(35, 13).Value = -2050.25  # M35: add -2050.25

$ws.Cells.Item(61, 8).Value = 1380  # H61: 1750 -> 1380
$ws.Cells.Item(61, 9).Value = 1380  # I61: 1750 -> 1380
$ws.Cells.Item(61, 11).Value = 1380  # K61: 1750 -> 1380
$ws.Cells.Item(61, 13).Value = -1168  # M61: -1538 -> -1168

$ws.Cells.Item(132, 8).Value = 35384.777  # H132: 33654 -> 35384.777
$ws.Cells.Item(132, 9).Value = 39183.25  # I132: 39264.5 -> 39183.25
$ws.Cells.Item(132, 10).Value = 4997  # J132: 3731.3333 -> 4997
$ws.Cells.Item(132, 11).Value = 117549.75  # K132: 117793.5 -> 117549.75
$ws.Cells.Item(132, 12).Value = 14991  # L132: 11193.9999 -> 14991
$ws.Cells.Item(132, 13).Value = -115019.75  # M132: -115263.5 -> -115019.75
$ws.Cells.Item(132, 14).Value = -20051  # N132: -16253.9999 -> -20051

$ws.Cells.Item(136, 8).Value = 1380  # H136: 1750 -> 1380
$ws.Cells.Item(136, 9).Value = 1380  # I136: 1750 -> 1380
$ws.Cells.Item(136, 11).Value = 4140  # K136: 5250 -> 4140
$ws.Cells.Item(136, 13).Value = -1590  # M136: -2700 -> -1590

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2058.7  # H134: 2118.6 -> 2058.7
$ws.Cells.Item(134, 9).Value = 2065.2222  # I134: 2131.7778 -> 2065.2222
$ws.Cells.Item(134, 11).Value = 6195.6666  # K134: 6395.3334 -> 6195.6666
$ws.Cells.Item(134, 13).Value = -3660.6666  # M134: -3860.3334 -> -3660.6666

$ws.Cells.Item(138, 8).Value = 84999.5  # H138: 85000 -> 84999.5
$ws.Cells.Item(138, 10).Value = 84999.5  # J138: 85000 -> 84999.5
$ws.Cells.Item(138, 12).Value = 84999.5  # L138: 85000 -> 84999.5
$ws.Cells.Item(138, 14).Value = -95279.5  # N138: -95280 -> -95279.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(55, 8).Value = 10000  # H55: 0 -> 10000
$ws.Cells.Item(55, 9).Value = 10000  # I55: 0 -> 10000
$ws.Cells.Item(55, 11).Value = 10000  # K55: 0 -> 10000
$ws.Cells.Item(55, 13).Value = -9685  # M55: add -9685

$ws.Cells.Item(99, 8).Value = 201619.6  # H99: 168391.33 -> 201619.6
$ws.Cells.Item(99, 10).Value = 0  # J99: 2250 -> 0
$ws.Cells.Item(99, 12).Value = 0  # L99: 2250 -> 0
$ws.Cells.Item(99, 14).ClearContents()  # N99: remove (was -5246)

$ws.Cells.Item(116, 8).Value = 424999  # H116: 599998.5 -> 424999
$ws.Cells.Item(116, 10).Value = 424999  # J116: 599998.5 -> 424999
$ws.Cells.Item(116, 12).Value = 424999  # L116: 599998.5 -> 424999
$ws.Cells.Item(116, 14).Value = -434177  # N116: -609176.5 -> -434177

$ws.Cells.Item(126, 8).Value = 201619.6  # H126: 168391.33 -> 201619.6
$ws.Cells.Item(126, 10).Value = 0  # J126: 2250 -> 0
$ws.Cells.Item(126, 12).Value = 0  # L126: 6750 -> 0
$ws.Cells.Item(126, 14).ClearContents()  # N126: remove (was -11690)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 1510  # H18: 946 -> 1510
$ws.Cells.Item(18, 9).Value = 1510  # I18: 946 -> 1510
$ws.Cells.Item(18, 11).Value = 4530  # K18: 2838 -> 4530
$ws.Cells.Item(18, 13).Value = -4361  # M18: -2669 -> -4361

$ws.Cells.Item(33, 8).Value = 580.8333  # H33: 629.0909 -> 580.8333
$ws.Cells.Item(33, 9).Value = 130  # I33: 140 -> 130
$ws.Cells.Item(33, 10).Value = 671  # J33: 678 -> 671
$ws.Cells.Item(33, 11).Value = 780  # K33: 840 -> 780
$ws.Cells.Item(33, 12).Value = 4026  # L33: 4068 -> 4026
$ws.Cells.Item(33, 13).Value = -497  # M33: -557 -> -497
$ws.Cells.Item(33, 14).Value = -4592  # N33: -4634 -> -4592

$ws.Cells.Item(68, 8).Value = 1992.3334  # H68: 2103.5 -> 1992.3334
$ws.Cells.Item(68, 10).Value = 2155.6667  # J68: 2366.2 -> 2155.6667
$ws.Cells.Item(68, 12).Value = 6467.000100000001  # L68: 7098.599999999999 -> 6467.000100000001
$ws.Cells.Item(68, 14).Value = -8089.000100000001  # N68: -8720.599999999999 -> -8089.000100000001

$ws.Cells.Item(71, 8).Value = 1992.3334  # H71: 2103.5 -> 1992.3334
$ws.Cells.Item(71, 10).Value = 2155.6667  # J71: 2366.2 -> 2155.6667
$ws.Cells.Item(71, 12).Value = 19401.0003  # L71: 21295.8 -> 19401.0003
$ws.Cells.Item(71, 14).Value = -27513.0003  # N71: -29407.8 -> -27513.0003

$ws.Cells.Item(107, 8).Value = 800.9677  # H107: 867.8570999999999 -> 800.9677
$ws.Cells.Item(107, 9).Value = 519.1111  # I107: 561.125 -> 519.1111
$ws.Cells.Item(107, 10).Value = 916.2727  # J107: 990.55 -> 916.2727
$ws.Cells.Item(107, 11).Value = 1557.3333  # K107: 1683.375 -> 1557.3333
$ws.Cells.Item(107, 12).Value = 2748.8181  # L107: 2971.65 -> 2748.8181
$ws.Cells.Item(107, 13).Value = 362.6667000000002  # M107: 236.625 -> 362.6667000000002
$ws.Cells.Item(107, 14).Value = -6588.8181  # N107: -6811.65 -> -6588.8181

$ws.Cells.Item(128, 8).Value = 344181.62  # H128: 344181.88 -> 344181.62
$ws.Cells.Item(128, 9).Value = 344181.62  # I128: 344181.88 -> 344181.62
$ws.Cells.Item(128, 11).Value = 1032544.86  # K128: 1032545.64 -> 1032544.86
$ws.Cells.Item(128, 13).Value = -1027564.86  # M128: -1027565.64 -> -1027564.86

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5667.684  # H70: 5734.25 -> 5667.684
$ws.Cells.Item(70, 9).Value = 5505  # I70: 5619.923 -> 5505
$ws.Cells.Item(70, 11).Value = 5505  # K70: 5619.923 -> 5505
$ws.Cells.Item(70, 13).Value = -5235  # M70: -5349.923 -> -5235

$ws.Cells.Item(73, 8).Value = 5667.684  # H73: 5734.25 -> 5667.684
$ws.Cells.Item(73, 9).Value = 5505  # I73: 5619.923 -> 5505
$ws.Cells.Item(73, 11).Value = 5505  # K73: 5619.923 -> 5505
$ws.Cells.Item(73, 13).Value = -4569  # M73: -4683.923 -> -4569

$ws.Cells.Item(80, 8).Value = 5212.857  # H80: 5199 -> 5212.857
$ws.Cells.Item(80, 9).Value = 4500  # I80: 5000 -> 4500
$ws.Cells.Item(80, 10).Value = 5747.5  # J80: 5331.6665 -> 5747.5
$ws.Cells.Item(80, 11).Value = 4500  # K80: 5000 -> 4500
$ws.Cells.Item(80, 12).Value = 5747.5  # L80: 5331.6665 -> 5747.5
$ws.Cells.Item(80, 13).Value = -3502  # M80: -4002 -> -3502
$ws.Cells.Item(80, 14).Value = -7743.5  # N80: -7327.6665 -> -7743.5

$ws.Cells.Item(83, 8).Value = 5212.857  # H83: 5199 -> 5212.857
$ws.Cells.Item(83, 9).Value = 4500  # I83: 5000 -> 4500
$ws.Cells.Item(83, 10).Value = 5747.5  # J83: 5331.6665 -> 5747.5
$ws.Cells.Item(83, 11).Value = 22500  # K83: 25000 -> 22500
$ws.Cells.Item(83, 12).Value = 28737.5  # L83: 26658.3325 -> 28737.5
$ws.Cells.Item(83, 13).Value = -17508  # M83: -20008 -> -17508
$ws.Cells.Item(83, 14).Value = -38721.5  # N83: -36642.3325 -> -38721.5

$ws.Cells.Item(102, 8).Value = 2274.5715  # H102: 2522.182 -> 2274.5715
$ws.Cells.Item(102, 9).Value = 1802.909  # I102: 1966.5 -> 1802.909
$ws.Cells.Item(102, 11).Value = 1802.909  # K102: 1966.5 -> 1802.909
$ws.Cells.Item(102, 13).Value = -180.9090000000001  # M102: -344.5 -> -180.9090000000001

$ws.Cells.Item(114, 8).Value = 150000  # H114: 0 -> 150000
$ws.Cells.Item(114, 10).Value = 150000  # J114: 0 -> 150000
$ws.Cells.Item(114, 12).Value = 150000  # L114: 0 -> 150000
$ws.Cells.Item(114, 14).Value = -158678  # N114: add -158678

$ws.Cells.Item(126, 8).Value = 3166.2856  # H126: 3181.2144 -> 3166.2856
$ws.Cells.Item(126, 9).Value = 3029.9092  # I126: 3048.9092 -> 3029.9092
$ws.Cells.Item(126, 11).Value = 9089.7276  # K126: 9146.7276 -> 9089.7276
$ws.Cells.Item(126, 13).Value = -6619.7276  # M126: -6676.7276 -> -6619.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6999.8184  # H7: 7569.9 -> 6999.8184
$ws.Cells.Item(7, 9).Value = 8042.7144  # I7: 8099.857 -> 8042.7144
$ws.Cells.Item(7, 10).Value = 5174.75  # J7: 6333.3335 -> 5174.75
$ws.Cells.Item(7, 11).Value = 8042.7144  # K7: 8099.857 -> 8042.7144
$ws.Cells.Item(7, 12).Value = 5174.75  # L7: 6333.3335 -> 5174.75
$ws.Cells.Item(7, 13).Value = -7930.7144  # M7: -7987.857 -> -7930.7144
$ws.Cells.Item(7, 14).Value = -5398.75  # N7: -6557.3335 -> -5398.75

$ws.Cells.Item(46, 8).Value = 2499  # H46: 2497.625 -> 2499
$ws.Cells.Item(46, 9).Value = 2499  # I46: 2497.625 -> 2499
$ws.Cells.Item(46, 11).Value = 2499  # K46: 2497.625 -> 2499
$ws.Cells.Item(46, 13).Value = -2311  # M46: -2309.625 -> -2311

$ws.Cells.Item(55, 8).Value = 446.27274  # H55: 464.2 -> 446.27274
$ws.Cells.Item(55, 9).Value = 470.125  # I55: 474 -> 470.125
$ws.Cells.Item(55, 10).Value = 382.66666  # J55: 425 -> 382.66666
$ws.Cells.Item(55, 11).Value = 470.125  # K55: 474 -> 470.125
$ws.Cells.Item(55, 12).Value = 382.66666  # L55: 425 -> 382.66666
$ws.Cells.Item(55, 13).Value = -297.125  # M55: -301 -> -297.125
$ws.Cells.Item(55, 14).Value = -728.66666  # N55: -771 -> -728.66666

$ws.Cells.Item(122, 8).Value = 12949.091  # H122: 13949.7 -> 12949.091
$ws.Cells.Item(122, 9).Value = 15930.25  # I122: 20000 -> 15930.25
$ws.Cells.Item(122, 10).Value = 4999.3335  # J122: 4874.25 -> 4999.3335
$ws.Cells.Item(122, 11).Value = 47790.75  # K122: 60000 -> 47790.75
$ws.Cells.Item(122, 12).Value = 14998.0005  # L122: 14622.75 -> 14998.0005
$ws.Cells.Item(122, 13).Value = -45340.75  # M122: -57550 -> -45340.75
$ws.Cells.Item(122, 14).Value = -19898.0005  # N122: -19522.75 -> -19898.0005

$ws.Cells.Item(126, 8).Value = 6999.8184  # H126: 7569.9 -> 6999.8184
$ws.Cells.Item(126, 9).Value = 8042.7144  # I126: 8099.857 -> 8042.7144
$ws.Cells.Item(126, 10).Value = 5174.75  # J126: 6333.3335 -> 5174.75
$ws.Cells.Item(126, 11).Value = 24128.1432  # K126: 24299.571 -> 24128.1432
$ws.Cells.Item(126, 12).Value = 15524.25  # L126: 19000.0005 -> 15524.25
$ws.Cells.Item(126, 13).Value = -21658.1432  # M126: -21829.571 -> -21658.1432
$ws.Cells.Item(126, 14).Value = -20464.25  # N126: -23940.0005 -> -20464.25

$ws.Cells.Item(132, 8).Value = 2797.4614  # H132: 2837.72 -> 2797.4614
$ws.Cells.Item(132, 9).Value = 2170  # I132: 2207.45 -> 2170
$ws.Cells.Item(132, 10).Value = 6248.5  # J132: 5358.8 -> 6248.5
$ws.Cells.Item(132, 11).Value = 6510  # K132: 6622.349999999999 -> 6510
$ws.Cells.Item(132, 12).Value = 18745.5  # L132: 16076.4 -> 18745.5
$ws.Cells.Item(132, 13).Value = -3980  # M132: -4092.349999999999 -> -3980
$ws.Cells.Item(132, 14).Value = -23805.5  # N132: -21136.4 -> -23805.5

$ws.Cells.Item(136, 8).Value = 2654.818  # H136: 2883.1 -> 2654.818
$ws.Cells.Item(136, 9).Value = 1627.4  # I136: 1820.5385 -> 1627.4
$ws.Cells.Item(136, 11).Value = 4882.200000000001  # K136: 5461.6155 -> 4882.200000000001
$ws.Cells.Item(136, 13).Value = -2332.200000000001  # M136: -2911.6155 -> -2332.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7392.8184  # H62: 7813.8335 -> 7392.8184
$ws.Cells.Item(62, 9).Value = 7791.1875  # I62: 8281.923000000001 -> 7791.1875
$ws.Cells.Item(62, 10).Value = 6330.5  # J62: 6596.8 -> 6330.5
$ws.Cells.Item(62, 11).Value = 7791.1875  # K62: 8281.923000000001 -> 7791.1875
$ws.Cells.Item(62, 12).Value = 6330.5  # L62: 6596.8 -> 6330.5
$ws.Cells.Item(62, 13).Value = -7167.1875  # M62: -7657.923000000001 -> -7167.1875
$ws.Cells.Item(62, 14).Value = -7578.5  # N62: -7844.8 -> -7578.5

$ws.Cells.Item(65, 8).Value = 7392.8184  # H65: 7813.8335 -> 7392.8184
$ws.Cells.Item(65, 9).Value = 7791.1875  # I65: 8281.923000000001 -> 7791.1875
$ws.Cells.Item(65, 10).Value = 6330.5  # J65: 6596.8 -> 6330.5
$ws.Cells.Item(65, 11).Value = 38955.9375  # K65: 41409.61500000001 -> 38955.9375
$ws.Cells.Item(65, 12).Value = 31652.5  # L65: 32984 -> 31652.5
$ws.Cells.Item(65, 13).Value = -35835.9375  # M65: -38289.61500000001 -> -35835.9375
$ws.Cells.Item(65, 14).Value = -37892.5  # N65: -39224 -> -37892.5

$ws.Cells.Item(100, 8).Value = 2234.2778  # H100: 2240.111 -> 2234.2778
$ws.Cells.Item(100, 9).Value = 2754  # I100: 2975.5833 -> 2754
$ws.Cells.Item(100, 10).Value = 883  # J100: 769.1667 -> 883
$ws.Cells.Item(100, 11).Value = 5508  # K100: 5951.1666 -> 5508
$ws.Cells.Item(100, 12).Value = 1766  # L100: 1538.3334 -> 1766
$ws.Cells.Item(100, 13).Value = -4967  # M100: -5410.1666 -> -4967
$ws.Cells.Item(100, 14).Value = -2848  # N100: -2620.3334 -> -2848

$ws.Cells.Item(122, 8).Value = 1922.1936  # H122: 2008.4828 -> 1922.1936
$ws.Cells.Item(122, 9).Value = 1455.5416  # I122: 1526.8636 -> 1455.5416
$ws.Cells.Item(122, 11).Value = 4366.6248  # K122: 4580.5908 -> 4366.6248
$ws.Cells.Item(122, 13).Value = -1916.6248  # M122: -2130.5908 -> -1916.6248

$ws.Cells.Item(126, 8).Value = 3575.1  # H126: 4450 -> 3575.1
$ws.Cells.Item(126, 9).Value = 3194.7778  # I126: 4025.3333 -> 3194.7778
$ws.Cells.Item(126, 11).Value = 9584.3334  # K126: 12075.9999 -> 9584.3334
$ws.Cells.Item(126, 13).Value = -7114.3334  # M126: -9605.999899999999 -> -7114.3334

$ws.Cells.Item(132, 8).Value = 13755.19  # H132: 15664.723 -> 13755.19
$ws.Cells.Item(132, 9).Value = 12776.214  # I132: 13690.23 -> 12776.214
$ws.Cells.Item(132, 10).Value = 15713.143  # J132: 20798.4 -> 15713.143
$ws.Cells.Item(132, 11).Value = 38328.642  # K132: 41070.69 -> 38328.642
$ws.Cells.Item(132, 12).Value = 47139.429  # L132: 62395.2 -> 47139.429
$ws.Cells.Item(132, 13).Value = -35798.642  # M132: -38540.69 -> -35798.642
$ws.Cells.Item(132, 14).Value = -52199.429  # N132: -67455.20000000001 -> -52199.429

$ws.Cells.Item(136, 8).Value = 939.7895  # H136: 944.8421 -> 939.7895
$ws.Cells.Item(136, 9).Value = 958.6667  # I136: 964 -> 958.6667
$ws.Cells.Item(136, 11).Value = 2876.0001  # K136: 2892 -> 2876.0001
$ws.Cells.Item(136, 13).Value = -326.0001000000002  # M136: -342 -> -326.0001000000002
